$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "updated" to "Tabelle1"
$ws.Name = "Tabelle1"

# Move/update the current selection to B10 (single cell)
$ws.Range("B10").Select()
